$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2017738.9
$ws.Range("I112").Value = 3110
$ws.Range("J112").Value = 2305543
$ws.Range("K112").Value = 9330
$ws.Range("L112").Value = 6916629
$ws.Range("M112").Value = -8222
$ws.Range("N112").Value = -6918845
$ws.Range("H141").Value = 1362.1052
$ws.Range("I141").Value = 1275.4706
$ws.Range("J141").Value = 2098.5
$ws.Range("K141").Value = 3826.4118
$ws.Range("L141").Value = 6295.5
$ws.Range("M141").Value = 1353.5882
$ws.Range("N141").Value = -16655.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16035587
$ws.Range("I32").Value = 16604886
$ws.Range("K32").Value = 16604886
$ws.Range("M32").Value = -16604599
$ws.Range("H45").Value = 4135.0835
$ws.Range("I45").Value = 3887.1
$ws.Range("J45").Value = 5375
$ws.Range("K45").Value = 3887.1
$ws.Range("L45").Value = 5375
$ws.Range("M45").Value = -3510.1
$ws.Range("N45").Value = -6129
$ws.Range("H55").Value = 14524
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 2115.7932
$ws.Range("I61").Value = 1978.0426
$ws.Range("J61").Value = 2704.3635
$ws.Range("K61").Value = 1978.0426
$ws.Range("L61").Value = 2704.3635
$ws.Range("M61").Value = -1766.0426
$ws.Range("N61").Value = -3128.3635
$ws.Range("H74").Value = 2427.9788
$ws.Range("I74").Value = 2400.4092
$ws.Range("K74").Value = 2400.4092
$ws.Range("M74").Value = -1526.4092
$ws.Range("H77").Value = 2427.9788
$ws.Range("I77").Value = 2400.4092
$ws.Range("K77").Value = 12002.046
$ws.Range("M77").Value = -7634.046
$ws.Range("H110").Value = 1545.8889
$ws.Range("I110").Value = 1301.6875
$ws.Range("K110").Value = 1301.6875
$ws.Range("M110").Value = 743.3125
$ws.Range("H132").Value = 2582.907
$ws.Range("I132").Value = 1970.4
$ws.Range("J132").Value = 3996.3845
$ws.Range("K132").Value = 5911.200000000001
$ws.Range("L132").Value = 11989.1535
$ws.Range("M132").Value = -3381.200000000001
$ws.Range("N132").Value = -17049.1535
$ws.Range("H136").Value = 2115.7932
$ws.Range("I136").Value = 1978.0426
$ws.Range("J136").Value = 2704.3635
$ws.Range("K136").Value = 5934.1278
$ws.Range("L136").Value = 8113.0905
$ws.Range("M136").Value = -3384.1278
$ws.Range("N136").Value = -13213.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H42").Value = 347005.5
$ws.Range("J42").Value = 347005.5
$ws.Range("L42").Value = 347005.5
$ws.Range("N42").Value = -347661.5
$ws.Range("H94").Value = 2818
$ws.Range("I94").Value = 2166.75
$ws.Range("K94").Value = 2166.75
$ws.Range("M94").Value = -1715.75
$ws.Range("H105").Value = 3562.25
$ws.Range("I105").Value = 3163.3333
$ws.Range("K105").Value = 3163.3333
$ws.Range("M105").Value = -1416.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1097.6666
$ws.Range("I16").Value = 928.9286
$ws.Range("J16").Value = 1688.25
$ws.Range("K16").Value = 928.9286
$ws.Range("L16").Value = 1688.25
$ws.Range("M16").Value = -641.9286
$ws.Range("N16").Value = -2262.25
$ws.Range("H31").Value = 2272.4707
$ws.Range("I31").Value = 2272.4707
$ws.Range("K31").Value = 2272.4707
$ws.Range("M31").Value = -1977.4707
$ws.Range("H34").Value = 2272.4707
$ws.Range("I34").Value = 2272.4707
$ws.Range("K34").Value = 2272.4707
$ws.Range("M34").Value = -2070.4707
$ws.Range("H62").Value = 1968.25
$ws.Range("I62").Value = 1624.6666
$ws.Range("K62").Value = 1624.6666
$ws.Range("M62").Value = -1000.6666
$ws.Range("H65").Value = 1968.25
$ws.Range("I65").Value = 1624.6666
$ws.Range("K65").Value = 8123.333000000001
$ws.Range("M65").Value = -5003.333000000001
$ws.Range("H105").Value = 2256.3572
$ws.Range("J105").Value = 2477.7778
$ws.Range("L105").Value = 2477.7778
$ws.Range("N105").Value = -5971.7778
$ws.Range("H122").Value = 3849522.5
$ws.Range("I122").Value = 4351201.5
$ws.Range("J122").Value = 3316.6667
$ws.Range("K122").Value = 13053604.5
$ws.Range("L122").Value = 9950.000100000001
$ws.Range("M122").Value = -13051154.5
$ws.Range("N122").Value = -14850.0001
$ws.Range("H132").Value = 2846.7144
$ws.Range("I132").Value = 2783.8823
$ws.Range("K132").Value = 8351.6469
$ws.Range("M132").Value = -5821.6469

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2671.75
$ws.Range("I18").Value = 475.6
$ws.Range("K18").Value = 1426.8
$ws.Range("M18").Value = -1257.8
$ws.Range("H23").Value = 376.9091
$ws.Range("I23").Value = 326.83334
$ws.Range("J23").Value = 437
$ws.Range("K23").Value = 980.5000200000001
$ws.Range("L23").Value = 1311
$ws.Range("M23").Value = -745.5000200000001
$ws.Range("N23").Value = -1781
$ws.Range("H92").Value = 1172.375
$ws.Range("H132").Value = 1111.25
$ws.Range("I132").Value = 1074
$ws.Range("J132").Value = 1148.5
$ws.Range("K132").Value = 9666
$ws.Range("L132").Value = 10336.5
$ws.Range("M132").Value = -7136
$ws.Range("N132").Value = -15396.5
$ws.Range("H137").Value = 2466.1924
$ws.Range("I137").Value = 712
$ws.Range("K137").Value = 2136
$ws.Range("M137").Value = 2964
$ws.Range("H139").Value = 2600.85
$ws.Range("I139").Value = 2260.6155
$ws.Range("J139").Value = 3232.7144
$ws.Range("K139").Value = 6781.8465
$ws.Range("L139").Value = 9698.143199999999
$ws.Range("M139").Value = -1641.8465
$ws.Range("N139").Value = -19978.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3666316.2
$ws.Range("I11").Value = 4656326
$ws.Range("K11").Value = 4656326
$ws.Range("M11").Value = -4656187
$ws.Range("H102").Value = 1796
$ws.Range("I102").Value = 1796
$ws.Range("K102").Value = 1796
$ws.Range("M102").Value = -174
$ws.Range("H132").Value = 3243
$ws.Range("I132").Value = 3321.6155
$ws.Range("K132").Value = 9964.8465
$ws.Range("M132").Value = -7434.8465
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 14500329
$ws.Range("I40").Value = 19614828
$ws.Range("J40").Value = 9249
$ws.Range("K40").Value = 19614828
$ws.Range("L40").Value = 9249
$ws.Range("M40").Value = -19614692
$ws.Range("N40").Value = -9521
$ws.Range("H46").Value = 4315.846
$ws.Range("I46").Value = 2895
$ws.Range("J46").Value = 4574.1816
$ws.Range("K46").Value = 2895
$ws.Range("L46").Value = 4574.1816
$ws.Range("M46").Value = -2707
$ws.Range("N46").Value = -4950.1816
$ws.Range("H82").Value = 16964.715
$ws.Range("I82").Value = 3009
$ws.Range("J82").Value = 30920.428
$ws.Range("K82").Value = 3009
$ws.Range("L82").Value = 30920.428
$ws.Range("M82").Value = -2648
$ws.Range("N82").Value = -31642.428
$ws.Range("H85").Value = 16964.715
$ws.Range("I85").Value = 3009
$ws.Range("J85").Value = 30920.428
$ws.Range("K85").Value = 3009
$ws.Range("L85").Value = 30920.428
$ws.Range("M85").Value = -1761
$ws.Range("N85").Value = -33416.428
$ws.Range("H100").Value = 4069.3333
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H132").Value = 84219.4
$ws.Range("I132").Value = 87436.875
$ws.Range("K132").Value = 262310.625
$ws.Range("M132").Value = -259780.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 95388.2
$ws.Range("J64").Value = 96735
$ws.Range("L64").Value = 96735
$ws.Range("N64").Value = -97231
$ws.Range("H67").Value = 95388.2
$ws.Range("J67").Value = 96735
$ws.Range("L67").Value = 96735
$ws.Range("N67").Value = -98451
$ws.Range("H81").Value = 2667.027
$ws.Range("I81").Value = 1826.3334
$ws.Range("J81").Value = 4219.077
$ws.Range("K81").Value = 3652.6668
$ws.Range("L81").Value = 8438.154
$ws.Range("M81").Value = -2591.6668
$ws.Range("N81").Value = -10560.154
$ws.Range("H84").Value = 2667.027
$ws.Range("I84").Value = 1826.3334
$ws.Range("J84").Value = 4219.077
$ws.Range("K84").Value = 18263.334
$ws.Range("L84").Value = 42190.77
$ws.Range("M84").Value = -12959.334
$ws.Range("N84").Value = -52798.77
$ws.Range("H107").Value = 444.91666
$ws.Range("I107").Value = 454.44446
$ws.Range("K107").Value = 1363.33338
$ws.Range("M107").Value = 556.66662
$ws.Range("H113").Value = 491.4762
$ws.Range("I113").Value = 453.73685
$ws.Range("K113").Value = 1361.21055
$ws.Range("M113").Value = 808.78945
$ws.Range("H132").Value = 2968.4783
$ws.Range("I132").Value = 3184.1538
$ws.Range("K132").Value = 9552.4614
$ws.Range("M132").Value = -7022.4614
$ws.Range("H136").Value = 28104.975
$ws.Range("I136").Value = 2078.423
$ws.Range("J136").Value = 80158.08
$ws.Range("K136").Value = 6235.268999999999
$ws.Range("L136").Value = 240474.24
$ws.Range("M136").Value = -3685.268999999999
$ws.Range("N136").Value = -245574.24
